$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 16.158065
$ws.Range("H2").Value = 48.474195
$ws.Range("I2").Value = 0.05027302757218571
$ws.Range("J2").Value = 0.0511740544428244
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.02264433333333333
$ws.Range("N2").Value = 0.06793299999999999
$ws.Range("O2").Value = 0.08454793804489194
$ws.Range("P2").Value = 0.08454793804489193
$ws.Range("Q2").Value = 0.3658886098816667
$ws.Range("R2").Value = 3.292997488935
$ws.Range("S2").Value = 0.004250480820502302
$ws.Range("T2").Value = 0.004326660784537844
$ws.Range("G3").Value = 16.158065
$ws.Range("H3").Value = 48.474195
$ws.Range("I3").Value = 0.05027302757218571
$ws.Range("J3").Value = 0.0511740544428244
$ws.Range("M3").Value = 0.245184
$ws.Range("N3").Value = 0.735552
$ws.Range("O3").Value = 0.915452061955108
$ws.Range("P3").Value = 0.9154520619551081
$ws.Range("Q3").Value = 3.96169900896
$ws.Range("R3").Value = 35.65529108064
$ws.Range("S3").Value = 0.04602254675168341
$ws.Range("T3").Value = 0.04684739365828656
$ws.Range("I4").Value = 0.8108637215090658
$ws.Range("J4").Value = 0.8253965642040216
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02264433333333333
$ws.Range("N4").Value = 0.06793299999999999
$ws.Range("O4").Value = 0.08454793804489194
$ws.Range("P4").Value = 0.08454793804489193
$ws.Range("Q4").Value = 5.901490604289222
$ws.Range("R4").Value = 53.113415438603
$ws.Range("S4").Value = 0.06855685568899901
$ws.Range("T4").Value = 0.06978557757278828
$ws.Range("I5").Value = 0.8108637215090658
$ws.Range("J5").Value = 0.8253965642040216
$ws.Range("M5").Value = 0.245184
$ws.Range("N5").Value = 0.735552
$ws.Range("O5").Value = 0.915452061955108
$ws.Range("P5").Value = 0.9154520619551081
$ws.Range("Q5").Value = 63.899036064448
$ws.Range("R5").Value = 575.0913245800321
$ws.Range("S5").Value = 0.7423068658200668
$ws.Range("T5").Value = 0.7556109866312334
$ws.Range("G6").Value = 10.026051
$ws.Range("H6").Value = 30.078153
$ws.Range("I6").Value = 0.03119432545686257
$ws.Range("J6").Value = 0.03175341105018045
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.02264433333333333
$ws.Range("N6").Value = 0.06793299999999999
$ws.Range("O6").Value = 0.08454793804489194
$ws.Range("P6").Value = 0.08454793804489193
$ws.Range("Q6").Value = 0.227033240861
$ws.Range("R6").Value = 2.043299167749
$ws.Range("S6").Value = 0.002637415896079012
$ws.Range("T6").Value = 0.002684685430184644
$ws.Range("G7").Value = 10.026051
$ws.Range("H7").Value = 30.078153
$ws.Range("I7").Value = 0.03119432545686257
$ws.Range("J7").Value = 0.03175341105018045
$ws.Range("M7").Value = 0.245184
$ws.Range("N7").Value = 0.735552
$ws.Range("O7").Value = 0.915452061955108
$ws.Range("P7").Value = 0.9154520619551081
$ws.Range("Q7").Value = 2.458227288384
$ws.Range("R7").Value = 22.124045595456
$ws.Range("S7").Value = 0.02855690956078356
$ws.Range("T7").Value = 0.02906872561999581
$ws.Range("G8").Value = 16.977099
$ws.Range("H8").Value = 33.954198
$ws.Range("I8").Value = 0.05282131035632833
$ws.Range("J8").Value = 0.03584533950516226
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.02264433333333333
$ws.Range("N8").Value = 0.06793299999999999
$ws.Range("O8").Value = 0.08454793804489194
$ws.Range("P8").Value = 0.08454793804489193
$ws.Range("Q8").Value = 0.384435088789
$ws.Range("R8").Value = 2.306610532734
$ws.Range("S8").Value = 0.004465932875456857
$ws.Range("T8").Value = 0.003030649543680576
$ws.Range("G9").Value = 16.977099
$ws.Range("H9").Value = 33.954198
$ws.Range("I9").Value = 0.05282131035632833
$ws.Range("J9").Value = 0.03584533950516226
$ws.Range("M9").Value = 0.245184
$ws.Range("N9").Value = 0.735552
$ws.Range("O9").Value = 0.915452061955108
$ws.Range("P9").Value = 0.9154520619551081
$ws.Range("Q9").Value = 4.162513041215999
$ws.Range("R9").Value = 24.975078247296
$ws.Range("S9").Value = 0.04835537748087147
$ws.Range("T9").Value = 0.03281468996148169
$ws.Range("G10").Value = 17.628366
$ws.Range("H10").Value = 52.885098
$ws.Range("I10").Value = 0.05484761510555757
$ws.Range("J10").Value = 0.05583063079781116
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.02264433333333333
$ws.Range("N10").Value = 0.06793299999999999
$ws.Range("O10").Value = 0.08454793804489194
$ws.Range("P10").Value = 0.08454793804489193
$ws.Range("Q10").Value = 0.399182595826
$ws.Range("R10").Value = 3.592643362434
$ws.Range("S10").Value = 0.004637252763854761
$ws.Range("T10").Value = 0.004720364713700573
$ws.Range("G11").Value = 17.628366
$ws.Range("H11").Value = 52.885098
$ws.Range("I11").Value = 0.05484761510555757
$ws.Range("J11").Value = 0.05583063079781116
$ws.Range("M11").Value = 0.245184
$ws.Range("N11").Value = 0.735552
$ws.Range("O11").Value = 0.915452061955108
$ws.Range("P11").Value = 0.9154520619551081
$ws.Range("Q11").Value = 4.322193289344
$ws.Range("R11").Value = 38.899739604096
$ws.Range("S11").Value = 0.05021036234170281
$ws.Range("T11").Value = 0.05111026608411059
